$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from B1 to C1, then set C1 value (new date column)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "13-01-2023"

$ws.Range("A2").Value = "1810 Renta variable"
$ws.Range("B2").Value = 127119.16
$ws.Range("C2").Value = 126964.74

$ws.Range("A3").Value = "1822 Raices Valores Negociables"
$ws.Range("B3").Value = 278044.48
$ws.Range("C3").Value = 278103.23

$ws.Range("A4").Value = "Adcap IOL Acciones Argentina"
$ws.Range("B4").Value = 58845.08
$ws.Range("C4").Value = 58867.61

$ws.Range("A5").Value = "Allaria Acciones"
$ws.Range("B5").Value = 39837.91
$ws.Range("C5").Value = 45298

$ws.Range("A6").Value = "Alpha Acciones"
$ws.Range("B6").Value = 108971.22
$ws.Range("C6").Value = 116254.12

$ws.Range("A7").Value = "Alpha Latam"
$ws.Range("B7").Value = 47.74
$ws.Range("C7").Value = 45.63

$ws.Range("A8").Value = "Alpha Mega"
$ws.Range("B8").Value = 156909.19
$ws.Range("C8").Value = 156874.24

$ws.Range("A9").Value = "Alpha Mercosur"
$ws.Range("B9").Value = 267484.69
$ws.Range("C9").Value = 282623.59

$ws.Range("A10").Value = "Alpha planeam equil"
$ws.Range("B10").Value = 6660.79
$ws.Range("C10").Value = 6658.4

$ws.Range("A11").Value = "Alpha renta balan global"
$ws.Range("B11").Value = 469823.83
$ws.Range("C11").Value = 470073.19

$ws.Range("A12").Value = "Argenfunds"
$ws.Range("B12").Value = 17157.68
$ws.Range("C12").Value = 17153.48

$ws.Range("A13").Value = "Arpenta ex Mercosur"
$ws.Range("B13").Value = 10593.73
$ws.Range("C13").Value = 10604.76

$ws.Range("A14").Value = "Balanz"
$ws.Range("B14").Value = 254254.97
$ws.Range("C14").Value = 260837.5

$ws.Range("A15").Value = "Bull Market"
$ws.Range("B15").Value = 79224.77
$ws.Range("C15").Value = 83969.24000000001

$ws.Range("A16").Value = "CMA acciones"
$ws.Range("B16").Value = 50132.93
$ws.Range("C16").Value = 91430.75999999999

$ws.Range("A17").Value = "Compass Crecimiento"
$ws.Range("B17").Value = 691572.96
$ws.Range("C17").Value = 691866.86

$ws.Range("A18").Value = "Consultatio Acciones Argentina"
$ws.Range("B18").Value = 145015.28
$ws.Range("C18").Value = 145634.23

$ws.Range("A19").Value = "Consultatio Renta Variable"
$ws.Range("B19").Value = 198285.29
$ws.Range("C19").Value = 198350.01

$ws.Range("A20").Value = "Delta Acciones"
$ws.Range("B20").Value = 74425.36
$ws.Range("C20").Value = 74377.03999999999

$ws.Range("A21").Value = "Delta Internacional"
$ws.Range("B21").Value = 1350.57
$ws.Range("C21").Value = 1350.43

$ws.Range("A22").Value = "Delta Latinoamerica"
$ws.Range("B22").Value = 3549.04
$ws.Range("C22").Value = 3546.75

$ws.Range("A23").Value = "Delta Select"
$ws.Range("B23").Value = 394376.61
$ws.Range("C23").Value = 394552.12

$ws.Range("A24").Value = "Delta gestion V"
$ws.Range("B24").Value = 162599.43
$ws.Range("C24").Value = 162879.73

$ws.Range("A25").Value = "FBA Acciones Argentinas"
$ws.Range("B25").Value = 158535.18
$ws.Range("C25").Value = 162107.15

$ws.Range("A26").Value = "FBA Calificado"
$ws.Range("B26").Value = 155943.39
$ws.Range("C26").Value = 158850.77

$ws.Range("A27").Value = "Fima Acciones"
$ws.Range("B27").Value = 291688.7
$ws.Range("C27").Value = 314277.38

$ws.Range("A28").Value = "Fima PB Acciones"
$ws.Range("B28").Value = 240254.59
$ws.Range("C28").Value = 252311.39

$ws.Range("A29").Value = "Gainvest Renta Variable"
$ws.Range("B29").Value = 328136.08
$ws.Range("C29").Value = 328105.35

$ws.Range("A30").Value = "Galileo Acciones"
$ws.Range("B30").Value = 2070784.18
$ws.Range("C30").Value = 2149672.19

$ws.Range("A31").Value = "Goal Acciones Argentinas"
$ws.Range("B31").Value = 33770.96
$ws.Range("C31").Value = 38760.08

$ws.Range("A32").Value = "Goal acciones plus"
$ws.Range("B32").Value = 6549.5
$ws.Range("C32").Value = 6539.17

$ws.Range("A33").Value = "HF Acciones Argentinas"
$ws.Range("B33").Value = 190537.56
$ws.Range("C33").Value = 190593.79

$ws.Range("A34").Value = "HF Acciones Lideres"
$ws.Range("B34").Value = 243337.28
$ws.Range("C34").Value = 243463.29

$ws.Range("A35").Value = "IAM Renta Variable"
$ws.Range("B35").Value = 55055.52
$ws.Range("C35").Value = 57695.35

$ws.Range("A36").Value = "IEB Value"
$ws.Range("B36").Value = 12858.09
$ws.Range("C36").Value = 12855.73

$ws.Range("A37").Value = "Lombardi"
$ws.Range("B37").Value = 36730.42
$ws.Range("C37").Value = 41466

$ws.Range("A38").Value = "MAF"
$ws.Range("B38").Value = 58335.01
$ws.Range("C38").Value = 58295.6

$ws.Range("A39").Value = "Megainver"
$ws.Range("B39").Value = 53903.51
$ws.Range("C39").Value = 53920.24

$ws.Range("A40").Value = "Pellegrini Acciones"
$ws.Range("B40").Value = 84374.39
$ws.Range("C40").Value = 84382.05

$ws.Range("A41").Value = "Premier Renta Variable"
$ws.Range("B41").Value = 62927.31
$ws.Range("C41").Value = 62929.25

$ws.Range("A42").Value = "Quinquela Acciones"
$ws.Range("B42").Value = 172138.29
$ws.Range("C42").Value = 172156.64

$ws.Range("A43").Value = "Rofex 20 Renta Variable"
$ws.Range("B43").Value = 110024.63
$ws.Range("C43").Value = 110045.92

$ws.Range("A44").Value = "SBS Acciones Argentina"
$ws.Range("B44").Value = 119468.33
$ws.Range("C44").Value = 148621.13

$ws.Range("A45").Value = "Schroeder RV"
$ws.Range("B45").Value = 669264.28
$ws.Range("C45").Value = 668821.3100000001

$ws.Range("A46").Value = "Supefondo RV"
$ws.Range("B46").Value = 514459.8
$ws.Range("C46").Value = 625236.7

$ws.Range("A47").Value = "Superfondo "
$ws.Range("B47").Value = 1098686.84
$ws.Range("C47").Value = 1099068.37

$ws.Range("A48").Value = "Toronto Trust Multimercado"
$ws.Range("B48").Value = 31412.83
$ws.Range("C48").Value = 38662.05

$ws.Range("A49").Value = "Toronto trust Argy"
$ws.Range("B49").Value = 15230.95
$ws.Range("C49").Value = 15184.49

$ws.Range("A50").Value = "avg"
$ws.Range("B50").Value = 216889.38
$ws.Range("C50").Value = 224423.69

$ws.Range("A51").Value = "total"
$ws.Range("B51").Value = 10410690.33
$ws.Range("C51").Value = 10772337.05
